$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right after the header (new row 2), shifting all
# existing data rows down by one (old row2 -> row3, old row3 -> row4, ...,
# old row77 -> row78).
$ws.Rows(2).Insert(-4121, 1)

# The freshly inserted row inherits header-ish formatting from the Insert
# call; strip that back to the plain/unstyled look the rest of the data
# rows use.
$ws.Range("A2:R2").ClearFormats()

# Column D carries a date number format (same as every other data row) -
# restore it by copying the format from the row right below (old row2,
# now row3).
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# New row 2 starts as a copy of the row beneath it (the former row 2),
# then a handful of cells get the genuinely new values from this week's
# entry.
$ws.Range("A2:R2").Value2 = $ws.Range("A3:R3").Value2

$ws.Range("D2").Value2 = 44631
$ws.Range("J2").Value2 = 120
$ws.Range("K2").Value2 = 29000
$ws.Range("L2").Value2 = 30000
$ws.Range("M2").Value2 = 29500
$ws.Range("P2").Value2 = 1180
